$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain decimal number (e.g. "590.21") would be
# auto-converted from text to a numeric type by Excel type-inference on .Value
# assignment. Force these specific cells to Text format first so the literal
# string (with its original formatting, trailing zeros, etc.) is preserved,
# matching how the source data is stored (inline string).
$textFormatCells = @("D5", "D6", "D10", "D14", "D18", "D19", "D20", "D22", "D26", "D27", "D28", "D29", "D33", "D35", "D36", "D37", "D39", "D40", "D41", "D42", "D45", "D47", "D49")
foreach ($addr in $textFormatCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values scraped by the GitHub Actions cron job.
$ws.Range('D2').Value = '68.846.28'
$ws.Range('E2').Value = '  +2.11%  '
$ws.Range('D3').Value = '3.305.14'
$ws.Range('E3').Value = '  +2.02%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = '590.21'
$ws.Range('E5').Value = '  +2.32%  '
$ws.Range('D6').Value = '186.87'
$ws.Range('E6').Value = '  +4.24%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('E8').Value = '  +1.26%  '
$ws.Range('E9').Value = '  +5.15%  '
$ws.Range('D10').Value = '6.74'
$ws.Range('E10').Value = '  -0.36%  '
$ws.Range('E11').Value = '  +2.88%  '
$ws.Range('D12').Value = '3.880.67'
$ws.Range('E12').Value = '  +2.16%  '
$ws.Range('E13').Value = '  +0.25%  '
$ws.Range('D14').Value = '29.07'
$ws.Range('E14').Value = '  +4.21%  '
$ws.Range('D15').Value = '68.853.55'
$ws.Range('E15').Value = '  +2.23%  '
$ws.Range('E16').Value = '  +3.88%  '
$ws.Range('D17').Value = '3.347.29'
$ws.Range('E17').Value = '  +3.32%  '
$ws.Range('D18').Value = '5.92'
$ws.Range('E18').Value = '  +2.09%  '
$ws.Range('D19').Value = '13.78'
$ws.Range('E19').Value = '  +3.25%  '
$ws.Range('D20').Value = '386.21'
$ws.Range('E20').Value = '  +3.40%  '
$ws.Range('E21').Value = '  +3.31%  '
$ws.Range('D22').Value = '71.76'
$ws.Range('E22').Value = '  +0.93%  '
$ws.Range('E23').Value = '  -0.23%  '
$ws.Range('E24').Value = '  +4.14%  '
$ws.Range('E25').Value = '  +2.53%  '
$ws.Range('D26').Value = '0.191'
$ws.Range('E26').Value = '  +5.47%  '
$ws.Range('D27').Value = '9.90'
$ws.Range('E27').Value = '  +2.80%  '
$ws.Range('D28').Value = '0.999'
$ws.Range('E28').Value = '  -0.42%  '
$ws.Range('D29').Value = '5.91'
$ws.Range('E29').Value = '  +5.67%  '
$ws.Range('E30').Value = '  +2.43%  '
$ws.Range('E31').Value = '  +5.53%  '
$ws.Range('E32').Value = '  +2.48%  '
$ws.Range('D33').Value = '7.26'
$ws.Range('E33').Value = '  +6.66%  '
$ws.Range('E34').Value = '  +0.03%  '
$ws.Range('D35').Value = '1.56'
$ws.Range('E35').Value = '  +4.52%  '
$ws.Range('D36').Value = '163.48'
$ws.Range('E36').Value = '  -0.39%  '
$ws.Range('D37').Value = '1.89'
$ws.Range('E37').Value = '  +2.59%  '
$ws.Range('E38').Value = '  -2.35%  '
$ws.Range('D39').Value = '26.96'
$ws.Range('E39').Value = '  +0.86%  '
$ws.Range('D40').Value = '6.79'
$ws.Range('E40').Value = '  -1.11%  '
$ws.Range('B41').Value = 'Filecoin'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D41').Value = '4.65'
$ws.Range('E41').Value = '  +5.63%  '
$ws.Range('B42').Value = 'dogwifhat'
$ws.Range('C42').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D42').Value = '2.67'
$ws.Range('E42').Value = '  +3.65%  '
$ws.Range('E43').Value = '  +0.58%  '
$ws.Range('E44').Value = '  +3.96%  '
$ws.Range('D45').Value = '41.46'
$ws.Range('E45').Value = '  +2.60%  '
$ws.Range('D46').Value = '2.650.71'
$ws.Range('E46').Value = '  -1.91%  '
$ws.Range('D47').Value = '342.83'
$ws.Range('E47').Value = '  -5.43%  '
$ws.Range('E48').Value = '  +3.43%  '
$ws.Range('D49').Value = '32.53'
$ws.Range('E49').Value = '  +6.56%  '
$ws.Range('E50').Value = '  +1.45%  '
$ws.Range('E51').Value = '  +0.53%  '
